$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 17 - spi flash chip (w25q16jvsniq, 2mbit)
$ws.Range("C17").Value = "2mbit"
$ws.Range("B17").Value = "w25q16jvsniq"
$ws.Range("D17").Value = "https://item.szlcsc.com/116655.html"

# Row 18 - 1v2 regulator (md5112)
$ws.Range("A18").Value = "1v2 reg"
$ws.Range("B18").Value = "md5112"
$ws.Range("D18").Value = "https://item.szlcsc.com/998676.html"

# Row 19 - lm317a
$ws.Range("B19").Value = "lm317a"

# Row 20 - xc6206p12
$ws.Range("B20").Value = "xc6206p12"
$ws.Range("D20").Value = "https://atta.szlcsc.com//upload/public/pdf/source/20130801/1457706628791.pdf"

# Row 21 - ld1117s12
$ws.Range("B21").Value = "ld1117s12"

# Row 22 - diode (1N5819WS)
$ws.Range("A22").Value = "diode"
$ws.Range("B22").Value = "1N5819WS"
$ws.Range("D22").Value = "https://atta.szlcsc.com//upload/public/pdf/source/20180614/C191023_3C6A6398B911F3A4C23E7538EE054643.pdf"

# Turn those cells into real hyperlinks (Datasheet column convention: display
# text equals the target URL), same as the rest of column D.
$ws.Hyperlinks.Add($ws.Range("D17"), "https://item.szlcsc.com/116655.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D18"), "https://item.szlcsc.com/998676.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D20"), "https://atta.szlcsc.com//upload/public/pdf/source/20130801/1457706628791.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D22"), "https://atta.szlcsc.com//upload/public/pdf/source/20180614/C191023_3C6A6398B911F3A4C23E7538EE054643.pdf") | Out-Null

# Copy the existing Datasheet-column hyperlink formatting (style) down onto the
# newly populated D cells, matching the existing look of D2:D16 (this must run
# AFTER Hyperlinks.Add, which otherwise mutates the cell's font flags).
$ws.Range("D16").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move selection to the new last cell, like Excel would leave it after typing.
$ws.Range("D22").Select()
